$d = $word.ActiveDocument

# 1. "The genetic modification..." paragraph: replace the opioids/vaccines/biotech
#    tail with the new "morphine, hydracodone, and many vaccines." ending.
$d.Content.Find.Execute(" Since then, researchers have been producing other useful medicines, such as opioids, and vaccines. This field also has an effect in biotechnology, where the modification of E. coli effects several other aspects that are not necessarily related to biology.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " Since then, researchers have been producing other useful medicines, such as morphine, hydracodone, and many vaccines. ", 2)

# 2. "Our research project..." paragraph: "Thebaine, a morphine precursor." -> "Thebaine."
$d.Content.Find.Execute(", a morphine precursor.", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)

# 3. Same paragraph: append the new "Manipulating 2 more genes..." sentence after
#    "get these results. "
$d.Content.Find.Execute("get these results. ", $true, $false, $false, $false, $false, $true, 1, $false, `
    "get these results. Manipulating 2 more genes from here allowed them for the creation of hydracodone.  ", 2)

# 4. "The model used for our code..." paragraph: remove the spell-check split
#    around "hydracodone" by re-writing the sentence (text unchanged).
$d.Content.Find.Execute("The model used for our code follows the engineering of E. coli to produce hydracodone pathways.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The model used for our code follows the engineering of E. coli to produce hydracodone pathways.", 2)

# 5. "One Next Slide..." paragraph: remove the spell-check split around "hydracodone".
$d.Content.Find.Execute("One Next Slide: This image displayed shows the pathway and how each of the chemicals links in order until hydracodone is produced.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "One Next Slide: This image displayed shows the pathway and how each of the chemicals links in order until hydracodone is produced.", 2)

# 6. "Next Slide: Our network..." paragraph: remove the spell-check split around
#    "hydracodone" (leave the "nodes is" grammar-check mark alone).
$d.Content.Find.Execute("Next Slide: Our network, which is the production model for hydracodone, has 14 ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Next Slide: Our network, which is the production model for hydracodone, has 14 ", 2)
